$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.088.93"
$ws.Range("D3").Value = "3.419.45"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'410.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'129.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").Value = "'0.640"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.26%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.740"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.87%  "
$ws.Range("D10").Value = "'0.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").Value = "'43.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.85%  "
$ws.Range("D12").Value = "'0.0000229"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +41.30%  "
$ws.Range("E13").Value = "  +9.00%  "
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "'21.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.13%  "
$ws.Range("D16").Value = "3.954.88"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "3.417.59"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").Value = "'12.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.38%  "
$ws.Range("E19").Value = "  +6.13%  "
$ws.Range("D20").Value = "62.022.47"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").Value = "'483.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +53.10%  "
$ws.Range("D22").Value = "'93.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.54%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "'13.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("E25").Value = "  +3.99%  "
$ws.Range("D26").Value = "'33.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.63%  "
$ws.Range("D27").Value = "'9.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.20%  "
$ws.Range("D28").Value = "'4.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "'7.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.54%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "'12.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.47%  "
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "'42.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.82%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'0.0507"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.75%  "
$ws.Range("D37").Value = "'53.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.98%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +7.94%  "
$ws.Range("D40").Value = "'3.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "'4.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.28%  "
$ws.Range("D44").Value = "'144.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("E45").Value = "  +16.38%  "
$ws.Range("D46").Value = "'2.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").Value = "'16.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("E48").Value = "  +19.64%  "
$ws.Range("D49").Value = "'22.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.41%  "
$ws.Range("D50").Value = "'2.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.90%  "
$ws.Range("D51").Value = "3.758.76"
$ws.Range("E51").Value = "  -0.80%  "
